# Insert a new data row at row 262 (pushing the existing rows 262-370 down
# to 263-371) and populate it with the new record's values. This mirrors
# the weekly refresh described in the commit message ("Fruta / hortaliza,
# semanal"): a new price observation is inserted near the top of the
# dataset and the rest of the history shifts down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 262..370 down to 263..371, carrying formatting along, and
# leave a blank row 262 ready to be filled in.
$ws.Rows.Item(262).Insert()

# Populate the newly inserted row 262 with the new observation.
$ws.Cells.Item(262, 1).Value = 7
$ws.Cells.Item(262, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(262, 3).Value = "Ñuble"
$ws.Cells.Item(262, 4).Value = 45229
$ws.Cells.Item(262, 5).Value = 16
$ws.Cells.Item(262, 6).Value = 100112032
$ws.Cells.Item(262, 7).Value = "Zapallo italiano"
$ws.Cells.Item(262, 8).Value = "Sin especificar"
$ws.Cells.Item(262, 9).Value = "Primera"
$ws.Cells.Item(262, 10).Value = 100
$ws.Cells.Item(262, 11).Value = 20000
$ws.Cells.Item(262, 12).Value = 20000
$ws.Cells.Item(262, 13).Value = 20000
$ws.Cells.Item(262, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(262, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(262, 16).Value = 400
$ws.Cells.Item(262, 17).Value = 50
$ws.Cells.Item(262, 18).Value = "Hortaliza"
